$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set E3 to a real boolean TRUE value instead of the string "True"
$ws.Range("E3").Value = $true
